$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("core i9 13900k")

# Append the new "6.5.0" benchmark row to the data table.
$ws.Range("A6").Value = "6.5.0"
$ws.Range("B6").Value = 35.334
$ws.Range("C6").Value = 25787196363
$ws.Range("D6").Value = 729812542

# Extend the linked chart series so they include the newly added row.
$chartObjs = $ws.ChartObjects()
$cols = @("B", "C", "D")
for ($i = 1; $i -le $chartObjs.Count; $i++) {
    $co = $chartObjs.Item($i)
    $chart = $co.Chart
    $seriesCollection = $chart.SeriesCollection()
    $series = $seriesCollection.Item(1)
    $col = $cols[$i - 1]
    $series.Formula = "=SERIES('core i9 13900k'!`$$col`$1,'core i9 13900k'!`$A`$2:`$A`$37,'core i9 13900k'!`$$col`$2:`$$col`$6,1)"
}

# Reflect the updated active selection on the sheet.
$ws.Activate()
$ws.Range("R22").Select()
